# Update the "contacts" sheet (sheet1) with new fields: nickname, position,
# department, category, status, phone - and fix the "MRS." -> "Mrs." title.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("contacts")
$ws2 = $wb.Worksheets.Item("signUp")

# Apply the same header highlight style used by A1:D1 to the new header cells
$ws1.Range("E1:J1").Interior.Color = $ws1.Range("A1").Interior.Color

# nickname column
$ws1.Range("E1").Value = "nickname"
$ws1.Range("E2").Value = "RS"
$ws1.Range("E3").Value = "VK"
$ws1.Range("E4").Value = "MR"

# position column
$ws1.Range("F1").Value = "position"
$ws1.Range("F2").Value = "Engineer"
$ws1.Range("F4").Value = "Crickter"
$ws1.Range("F3").Value = "Doctor"

# department column
$ws1.Range("G1").Value = "department"
$ws1.Range("G2").Value = "R&D"
$ws1.Range("G4").Value = "Womens Cricket"
$ws1.Range("G3").Value = "Surgeon"

# Fix title typo for Mithali Raj
$ws1.Range("A4").Value = "Mrs."

# category column
$ws1.Range("H1").Value = "category"
$ws1.Range("H3").Value = "Lead"
$ws1.Range("H2").Value = "Friend"
$ws1.Range("H4").Value = "Friend"

# status column
$ws1.Range("I1").Value = "status"
$ws1.Range("I2").Value = "Active"
$ws1.Range("I3").Value = "New"
$ws1.Range("I4").Value = "Active"

# phone column
$ws1.Range("J1").Value = "phone"
$ws1.Range("J2").Value = 11111
$ws1.Range("J3").Value = 22222
$ws1.Range("J4").Value = 33333

# Column widths (auto-fit best-fit columns) to mirror the commit's layout tweaks
$ws1.Columns.Item(5).AutoFit()
$ws1.Columns.Item(7).AutoFit()
$ws1.Columns.Item(8).AutoFit()
$ws1.Columns.Item(10).AutoFit()

# Make the contacts sheet portrait-oriented when printed
[void]($ws1.PageSetup.Orientation = 1)

# Switch active sheet/selection back to "contacts" (it was on "signUp")
$ws1.Activate() | Out-Null
$ws1.Range("J7").Select() | Out-Null

Write-Output "done"
